$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
    <w:p>
      <w:r>
        <w:t>Part A:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>The update transaction was still active and had not committed. The record was updated in memory but had not flushed to disk yet.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> The record in memory can be accesses and will show the updated result. When the transaction is terminated then the update is lost. When the system comes back up, after logging back in, the record is pulled from disk which has the value prior to the update. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Part B:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">The update transaction was still active and had not committed. Again, the record was updated in memory and not flushed to disk. The savepoint was set </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">but is deleted when a transaction terminates. Savepoints are only used to rollback the transaction that created the savepoint. A savepoint is a point in the transaction that can be </w:t>
      </w:r>
      <w:r>
        <w:t>referenced</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> using a rollback.  </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Part C:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>After the rollback to s1 all of the votes in the movie record were returned to their values of savepoint s1</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">The rollback to s2 rendered an error for being invalid. It is not possible </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">to rollback to a save point that </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>part  e step 11;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>I receive and error ORA-0054;resource busy and acquire with NOWAIT specified or timeout expired. This is expected becuase The nowaite clause is set so my trasaction will not wait for resources that are locked by another trasaction.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>part D</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Ask your partner to re-query the score of ' 'My Cousin Vinny' '. Does your partner see the change by your update ?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Recall that 2 phase locking schemes with read/write locks. Does Oracle follow the 2 phase locking schemes? Which kind of locking does Oracle follow? Explain your observation. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Next, we examine the effect of a change of transaction isolation level. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>You and your partner re-login to SQL*Plus and ensure that AUTOCOMMIT is OFF.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>--step 9 he cancels then requries and still sees 7.5 he can read but I have not commited the changes yet.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>-- I commit then he sees the correct score of 10.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>--IT releases the locks once the original process commit. When a process commits then it will release all of the locks that it has obtained.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>step 16;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>after I commit my partner still sees the score as 7.5; This is because the serializable trasaction level my partner's transaction will read the score as 7.5 until his trasaction commits. I can update the score of My Cousin Vinny as much as i want be he will not be able to see any other score.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>read commited snapshot reads the score of the last commited transation.  once my partner commits his trasaction then he will be able to see the updated score.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>update movie set score = (select max(score) from movie) where title = 'My Cousin Vinny';</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve"> releases the locks once the original process commit. When a process commits then it will release all of the locks that it has obtained.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Content.InsertXML($xml)
